$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Taxonsorteringsordning updated
$ws.Range("B2").Value = 90818

# Row 3: Taxonsorteringsordning updated
$ws.Range("B3").Value = 88637

# Rows 4 and 5: the two observation records (species data) are swapped,
# and the Taxonsorteringsordning (column B) values are refreshed.

# Row 4 becomes the former row 5 record, with new B value
$ws.Range("A4").Value = 112473083
$ws.Range("B4").Value = 89820
$ws.Range("D4").Value = "EN"
$ws.Range("E4").Value = 71
$ws.Range("F4").Value = "Urskogsporing"
$ws.Range("G4").Value = "Neoantrodia infirma"
$ws.Range("H4").Value = "(Renvall & Niemelä) Audet"
$ws.Range("Q4").Value = 518039
$ws.Range("R4").Value = 6790377
$ws.Range("Z4").Value = "13:17"
$ws.Range("AB4").Value = "13:17"

# Row 5 becomes the former row 4 record, with new B value
$ws.Range("A5").Value = 112472885
$ws.Range("B5").Value = 90818
$ws.Range("D5").Value = "VU"
$ws.Range("E5").Value = 4365
$ws.Range("F5").Value = "Smalfotad taggsvamp"
$ws.Range("G5").Value = "Hydnellum gracilipes"
$ws.Range("H5").Value = "(P.Karst) P.Karst"
$ws.Range("Q5").Value = 517956
$ws.Range("R5").Value = 6790407
$ws.Range("Z5").Value = "13:12"
$ws.Range("AB5").Value = "13:12"
